# Append a new day (2025-04-08) to every price sheet in the workbook,
# carrying forward the previous day's (row 37) price into the new row 38.

$wb = $excel.ActiveWorkbook
$newDate = "2025-04-08"

foreach ($ws in $wb.Worksheets) {
    # Price carried forward from the last existing row (row 37, column B).
    $lastPrice = $ws.Range("B37").Text

    # Use text format so values are written back as text (matching the
    # existing inlineStr/text cells in the sheet) rather than being
    # reinterpreted as dates or numbers.
    $ws.Range("A38").NumberFormat = "@"
    $ws.Range("A38").Value = $newDate
    $ws.Range("A38").Style = "Normal"

    $ws.Range("B38").NumberFormat = "@"
    $ws.Range("B38").Value = $lastPrice
    $ws.Range("B38").Style = "Normal"
}
